# Append 45 new data rows (102-146) to the worksheet, mirroring the
# existing repeating pattern of regcntr_id / machine_id already present
# in the sheet, continuing the device_id sequence, and keeping the same
# constant values for lang_code / is_active / cr_by / cr_dtimes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aCycle = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)
$bCycle = @(10021, 10022, 10023, 10024, 10025, 10026, 10027, 10028, 10029)

$startRow = 102
$startDeviceId = 3000121
$count = 45

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $idx = $i % 9

    $ws.Cells.Item($row, 1).Value = $aCycle[$idx]
    $ws.Cells.Item($row, 2).Value = $bCycle[$idx]
    $ws.Cells.Item($row, 3).Value = $startDeviceId + $i
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

$lastRow = $startRow + $count - 1

# Mirror the post-paste selection/view state captured in the diff: the
# newly added block is selected and the view is scrolled so row 129 is
# the top-left visible row.
$ws.Range("A102:G$lastRow").Select()
$excel.ActiveWindow.ScrollRow = 129

# Print setup metadata recorded on the sheet after the edit (portrait
# orientation page setup materializes the <pageSetup> element).
$ws.PageSetup.Orientation = 1
